$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 526 (A_SERIES_WITHOUTDIGI) needs to move up to row 510, with rows
# 510-525 shifting down by one row to become rows 511-526.
# Capture the values first (read top-to-bottom), then rewrite bottom-to-top
# so we never overwrite data before it has been read.

$firstRow = 510
$lastRow = 526

$values = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $values[$r] = @(
        $ws.Cells.Item($r, 1).Value2,
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2
    )
}

# New row (firstRow) gets what used to be in the last row.
$ws.Cells.Item($firstRow, 1).Value = $values[$lastRow][0]
$ws.Cells.Item($firstRow, 2).Value = $values[$lastRow][1]
$ws.Cells.Item($firstRow, 3).Value = $values[$lastRow][2]
$ws.Cells.Item($firstRow, 4).Value = $values[$lastRow][3]

# Every other row shifts down by one (was r-1, now r).
for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $src = $values[$r - 1]
    $ws.Cells.Item($r, 1).Value = $src[0]
    $ws.Cells.Item($r, 2).Value = $src[1]
    $ws.Cells.Item($r, 3).Value = $src[2]
    $ws.Cells.Item($r, 4).Value = $src[3]
}
